# ============================================================================
# Applies the "additional scraping" edit:
#   1. Renames the existing "ODI Batting" sheet's role: the workbook gains a
#      new first sheet "Player Info", the batting sheet becomes the 2nd
#      sheet, and a brand-new 3rd sheet "ODI Batting Extra" is appended.
#   2. "Player Info" gets a small 1-row table describing the player.
#   3. "ODI Batting" header D1 is renamed MATCH_CARD_LINK -> MATCH_CODE and
#      every row's D value (a full scorecard URL) is replaced by just the
#      trailing numeric match code. Also the two stray empty B38/B39 cells
#      are removed.
#   4. "ODI Batting Extra" is populated with extra per-match stats.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Helper: write a value into a cell while preserving its intended OOXML
# type. Values that "look numeric" (plain integers, decimals or
# percentages) get auto-converted to real numbers by Excel's normal
# assignment semantics, but several of the source cells must stay text
# (e.g. match codes, percentage strings). Forcing the NumberFormat to
# "@" (Text) before the assignment keeps the value as a string; clearing
# the format afterwards drops the now-unneeded text number format so we
# don't leave stray formatting behind.
# ----------------------------------------------------------------------
function Set-TextValue {
    param($cell, [string]$value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

function Set-NumberValue {
    param($cell, $value)
    $cell.Value = $value
}

function Style-Header {
    param($range)
    $range.Font.Bold = $true
    $range.Borders.LineStyle = 1
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4160
}

# Reproduce the page margins used throughout this workbook (0.75in sides,
# 1in top/bottom, 0.5in header/footer == 54/54/72/72/36/36 points) on a
# freshly-created sheet, which otherwise defaults to Excel's stock
# 0.7/0.75/0.3in margins.
function Set-StandardMargins {
    param($sheet)
    $sheet.PageSetup.LeftMargin = 54
    $sheet.PageSetup.RightMargin = 54
    $sheet.PageSetup.TopMargin = 72
    $sheet.PageSetup.BottomMargin = 72
    $sheet.PageSetup.HeaderMargin = 36
    $sheet.PageSetup.FooterMargin = 36
}

# ----------------------------------------------------------------------
# 1. Re-arrange sheets: rename the current (only) sheet out of the way,
#    insert the new "ODI Batting" sheet after it, then the new
#    "ODI Batting Extra" sheet after that, then rename the original
#    sheet to "Player Info" and the new batting sheet to "ODI Batting".
#    This ordering of operations reproduces the target sheetId
#    allocation (Player Info=1, ODI Batting=2, ODI Batting Extra=3).
# ----------------------------------------------------------------------
$battingSrc = $wb.Worksheets.Item(1)
$battingSrc.Name = "Player Info__tmp"

$battingNew = $wb.Worksheets.Add($null, $battingSrc)
$battingNew.Name = "ODI Batting__tmp"

$extra = $wb.Worksheets.Add($null, $battingNew)
$extra.Name = "ODI Batting Extra"

# Move all the existing batting data onto the new "ODI Batting" sheet.
$battingSrc.UsedRange.Copy($battingNew.Range("A1"))

# Finish renaming.
$battingSrc.Name = "Player Info"
$battingNew.Name = "ODI Batting"

$playerInfo = $battingSrc
$batting = $battingNew

Set-StandardMargins $batting
Set-StandardMargins $extra

# ----------------------------------------------------------------------
# 2. Clear out the old batting data from the "Player Info" sheet and
#    replace it with the small player-description table.
# ----------------------------------------------------------------------
$playerInfo.UsedRange.Clear()

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 0; $c -lt $piHeaders.Length; $c++) {
    Set-TextValue $playerInfo.Cells.Item(1, $c + 1) $piHeaders[$c]
}
Style-Header $playerInfo.Range("A1:D1")

$piRow = @("3720", "Kieran Omar Akeem Powell", "Left Handed", "Right Arm Medium")
for ($c = 0; $c -lt $piRow.Length; $c++) {
    Set-TextValue $playerInfo.Cells.Item(2, $c + 1) $piRow[$c]
}

# ----------------------------------------------------------------------
# 3. Fix up the "ODI Batting" sheet: rename the MATCH_CARD_LINK header to
#    MATCH_CODE, replace every URL in column D with just the trailing
#    match code, and drop the two stray empty B38/B39 cells.
# ----------------------------------------------------------------------
Set-TextValue $batting.Cells.Item(1, 4) "MATCH_CODE"

$lastRow = $batting.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $url = [string]$cell.Value()
    if ($url -match 'MatchCode=(\d+)') {
        Set-TextValue $cell $matches[1]
    }
}

$batting.Cells.Item(38, 2).ClearContents()
$batting.Cells.Item(39, 2).ClearContents()

# ----------------------------------------------------------------------
# 4. Populate "ODI Batting Extra".
# ----------------------------------------------------------------------
$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 0; $c -lt $extraHeaders.Length; $c++) {
    Set-TextValue $extra.Cells.Item(1, $c + 1) $extraHeaders[$c]
}
Style-Header $extra.Range("A1:F1")

# Each entry: MatchCode, BattingPosition(number or $null), NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("3625", 2, "2", "0", "10.06%", "NO"),
    @("3629", 1, "0", "0", "0.36%", "NO"),
    @("4001", 3, "0", "0", "0.40%", "NO"),
    @("4004", 3, "2", "0", "4.00%", "NO"),
    @("4005", 1, "0", "1", "4.23%", "NO"),
    @("4017", 3, "5", "0", "19.74%", "NO"),
    @("4018", 4, "0", "0", "5.29%", "NO"),
    @("4019", $null, $null, $null, $null, "NO"),
    @("4040", 2, "0", "0", "1.34%", "NO"),
    @("4043", 2, "1", "0", "8.70%", "NO"),
    @("4046", $null, $null, $null, $null, "NO"),
    @("4051", $null, $null, $null, $null, "NO"),
    @("4052", 1, "0", "0", $null, "NO"),
    @("4181", 5, "0", "0", "1.41%", "NO"),
    @("4213", $null, $null, $null, $null, "NO"),
    @("4216", 1, "3", "0", "5.61%", "NO"),
    @("4219", 1, "2", "1", "7.42%", "NO"),
    @("4220", 2, "1", "0", "2.61%", "NO"),
    @("4221", 1, "0", "0", $null, "NO"),
    @("4228", 1, "1", "0", "5.13%", "NO")
)

for ($i = 0; $i -lt $extraRows.Length; $i++) {
    $row = $extraRows[$i]
    $r = $i + 2

    Set-TextValue $extra.Cells.Item($r, 1) $row[0]

    if ($null -ne $row[1]) {
        Set-NumberValue $extra.Cells.Item($r, 2) $row[1]
    }
    if ($null -ne $row[2]) {
        Set-TextValue $extra.Cells.Item($r, 3) $row[2]
    }
    if ($null -ne $row[3]) {
        Set-TextValue $extra.Cells.Item($r, 4) $row[3]
    }
    if ($null -ne $row[4]) {
        Set-TextValue $extra.Cells.Item($r, 5) $row[4]
    }
    Set-TextValue $extra.Cells.Item($r, 6) $row[5]
}

# ----------------------------------------------------------------------
# Make "Player Info" the active sheet/tab, matching the original
# workbook-level activeTab="0".
# ----------------------------------------------------------------------
$playerInfo.Range("A1").Select()
$playerInfo.Activate()
